# Working Version with One ASIN
# Updates the "Forecast Comparison" sheet's forecast numbers and the
# "Summary" sheet's headline metrics to reflect the new (single-ASIN)
# forecast run.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: "Forecast Comparison" - Prophet / Amazon forecast columns
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Forecast Comparison")

# row, Prophet Forecast (B), Amazon Mean (C), P70 (D), P80 (E), P90 (F)
$forecastRows = @(
    @(2,  468.1561649076668, 701, 857, 1079, 1445),
    @(3,  532.1156751534114, 266, 321,  383,  480),
    @(4,  457.7327641261022, 139, 167,  196,  242),
    @(5,  317.5540424429,    173, 208,  247,  307),
    @(6,  215.533517665481,  121, 147,  175,  220),
    @(7,  195.3814903846867, 107, 130,  156,  197),
    @(8,  221.8929296655804, 113, 137,  165,  211),
    @(9,  240.1185641152562, 101, 122,  143,  176),
    @(10, 235.6136802083513, 101, 122,  144,  179),
    @(11, 232.2287374602311, 105, 128,  153,  193),
    @(12, 245.9224817702768, 142, 172,  207,  262),
    @(13, 260.7051362051131, 135, 165,  199,  254),
    @(14, 253.9997037880444, 137, 166,  198,  250),
    @(15, 230.8431160633121, 130, 158,  189,  238),
    @(16, 217.6013250812106, 131, 160,  193,  245),
    @(17, 223.3925503230307, 134, 163,  196,  250),
    @(18, 222.8242736541926, 124, 150,  180,  227),
    @(19, 186.8141746494535, 130, 158,  190,  240),
    @(20, 124.3343996311871, 140, 170,  207,  266),
    @(21, 80.51744823513057, 124, 151,  181,  228)
)

foreach ($entry in $forecastRows) {
    $r = $entry[0]
    $ws1.Cells.Item($r, 2).Value = $entry[1]
    $ws1.Cells.Item($r, 3).Value = $entry[2]
    $ws1.Cells.Item($r, 4).Value = $entry[3]
    $ws1.Cells.Item($r, 5).Value = $entry[4]
    $ws1.Cells.Item($r, 6).Value = $entry[5]
}

# ---------------------------------------------------------------------
# Sheet 2: "Summary" - headline metrics (column B), all stored as text
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Summary")

# row, new text value
$summaryRows = @(
    @(3,  "8"),
    @(4,  "666"),
    @(5,  "188"),
    @(6,  "170"),
    @(7,  "100"),
    @(8,  "18802 units"),
    @(9,  "4549"),
    @(10, "2648"),
    @(11, "1776"),
    @(12, "532"),
    @(13, "2024-12-01"),
    @(14, "81")
)

foreach ($entry in $summaryRows) {
    $r = $entry[0]
    $cell = $ws2.Cells.Item($r, 2)
    # Prefix with an apostrophe so Excel stores the value as literal text
    # instead of auto-converting number-/date-looking strings, then reset
    # the cell style so the resulting "number stored as text" quote-prefix
    # formatting isn't left behind on the cell.
    $cell.Value = "'" + $entry[1]
    $cell.Style = "Normal"
}
